# "Generate Report for Handback"
#
# For the a.md / b.md rows on each language sheet (zh-cn, de-de), mark the
# handback as complete:
#   - Status               -> "Handed back: in sync with en-US"
#   - Latest Target File   -> same file name as "Source File Name" (hyperlinked)
#   - Latest Handback File -> same file name as "Latest Handoff File" (hyperlinked)
#   - Latest Handback DateTime -> timestamp of the handback

function Get-HyperlinkAddress($ws, $cellAddr) {
    foreach ($h in $ws.Hyperlinks) {
        $hAddr = $h.Range.Address()
        if ($hAddr -eq $cellAddr) {
            return $h.Address()
        }
    }
    return ""
}

function Set-HandbackRow($ws, $row, $handbackStatus, $handbackDateTime, $srcRow) {
    # Latest Target File / Latest Handback File mirror the handoff info of
    # $srcRow (the row that was actually handed off - row 2 / "a.md" for
    # this fixture), regardless of which row is being marked handed-back.
    $srcAddr = "A" + $srcRow
    $handoffAddr = "C" + $srcRow
    $targetAddr = "E" + $row
    $handbackAddr = "F" + $row
    $dateAddr = "G" + $row

    $srcText = $ws.Range($srcAddr).Value2()
    $handoffText = $ws.Range($handoffAddr).Value2()

    $srcCellAddr = "`$A`$" + $srcRow
    $handoffCellAddr = "`$C`$" + $srcRow
    $srcUrl = Get-HyperlinkAddress $ws $srcCellAddr
    $handoffUrl = Get-HyperlinkAddress $ws $handoffCellAddr

    # Status: handed back, in sync with en-US
    $ws.Range("B" + $row).Value = $handbackStatus

    # Latest Target File = Source File Name (hyperlinked the same way as column A)
    $targetRange = $ws.Range($targetAddr)
    $ws.Hyperlinks.Add($targetRange, $srcUrl, "", "", $srcText)

    # Latest Handback File = Latest Handoff File (hyperlinked the same way as column C)
    $handbackRange = $ws.Range($handbackAddr)
    $ws.Hyperlinks.Add($handbackRange, $handoffUrl, "", "", $handoffText)

    # Latest Handback DateTime
    $ws.Range($dateAddr).Value = $handbackDateTime
}

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn 2 "Handed back: in sync with en-US" "2016-02-22 13:45:48" 2
Set-HandbackRow $wsZhCn 3 "Handed back: in sync with en-US" "2016-02-22 13:45:48" 2

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe 2 "Handed back: in sync with en-US" "2016-02-22 13:46:14" 2
Set-HandbackRow $wsDeDe 3 "Handed back: in sync with en-US" "2016-02-22 13:46:14" 2
